# Commit: "Remove policies from BAU files"
#
# On the BBNPPTY sheet, the "hard coal" row (row 2) and "lignite" row
# (row 14) allowed new builds of those technologies starting in 2028
# (columns I:AE, years 2028-2050, were flagged 1/true). This change
# removes that allowance by flipping those flags to 0/false, consistent
# with the rest of the fossil-fuel rows on the sheet which already ban
# new builds across the whole horizon.

$wb = $excel.ActiveWorkbook
$aboutSheet = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("BBNPPTY")

# hard coal (row 2) and lignite (row 14): disallow new builds 2028-2050
$ws.Range("I2:AE2").Value = 0
$ws.Range("I14:AE14").Value = 0

# Reflect the updated selection/scroll position on the BBNPPTY sheet,
# then restore "About" as the active sheet/tab (matches original file).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 12
$ws.Range("H14:AE14").Select()

$aboutSheet.Activate()
